$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily Update 키워드 10개 - apply scraped price/listing changes to rows

# Row 4
$ws.Range("E4").Value = '''46000'
# Row 10
$ws.Range("B10").Value = '큐디스 USB LED 스탠드'
$ws.Range("C10").Value = 'https://search.shopping.naver.com/gate.nhn?id=17060899687'
$ws.Range("D10").Value = 'https://shopping-phinf.pstatic.net/main_1706089/17060899687.20210917144350.jpg'
$ws.Range("E10").Value = '''4900'
$ws.Range("I10").Value = '큐디스'
$ws.Range("J10").Value = '큐디스'
# Row 11
$ws.Range("B11").Value = '듀플렉스 DP-310LS'
$ws.Range("C11").Value = 'https://search.shopping.naver.com/gate.nhn?id=9432703889'
$ws.Range("D11").Value = 'https://shopping-phinf.pstatic.net/main_9432703/9432703889.20201008151032.jpg'
$ws.Range("E11").Value = '''24750'
$ws.Range("I11").Value = '듀플렉스'
$ws.Range("J11").Value = '듀플렉스'
# Row 23
$ws.Range("E23").Value = '''23360'
# Row 26
$ws.Range("B26").Value = '이지넷유비쿼터스 넥스트 NEXT-211LAMP-W'
$ws.Range("C26").Value = 'https://search.shopping.naver.com/gate.nhn?id=27132503522'
$ws.Range("D26").Value = 'https://shopping-phinf.pstatic.net/main_2713250/27132503522.20210513084051.jpg'
$ws.Range("E26").Value = '''19900'
$ws.Range("I26").Value = '넥스트'
$ws.Range("J26").Value = '이지넷유비쿼터스'
# Row 27
$ws.Range("B27").Value = '브리츠 BE-LED10W'
$ws.Range("C27").Value = 'https://search.shopping.naver.com/gate.nhn?id=22265886605'
$ws.Range("D27").Value = 'https://shopping-phinf.pstatic.net/main_2226588/22265886605.20200323103143.jpg'
$ws.Range("E27").Value = '''39800'
$ws.Range("I27").Value = '브리츠'
$ws.Range("J27").Value = '브리츠'
# Row 28
$ws.Range("B28").Value = '플랜룩스 모티아이 LED 스탠드'
$ws.Range("C28").Value = 'https://search.shopping.naver.com/gate.nhn?id=21100161497'
$ws.Range("D28").Value = 'https://shopping-phinf.pstatic.net/main_2110016/21100161497.20211111154138.jpg'
$ws.Range("E28").Value = '''77800'
$ws.Range("I28").Value = '플랜룩스'
$ws.Range("J28").Value = '플랜룩스'
# Row 29
$ws.Range("B29").Value = '이지넷유비쿼터스 넥스트 NEXT-122LAMP-WC'
$ws.Range("C29").Value = 'https://search.shopping.naver.com/gate.nhn?id=21396473240'
$ws.Range("D29").Value = 'https://shopping-phinf.pstatic.net/main_2139647/21396473240.20211215184908.jpg'
$ws.Range("E29").Value = '''18900'
# Row 37
$ws.Range("E37").Value = '''16270'
# Row 52
$ws.Range("E52").Value = '''44680'
# Row 54
$ws.Range("B54").Value = '아이클 WJK-151C'
$ws.Range("C54").Value = 'https://search.shopping.naver.com/gate.nhn?id=6512970130'
$ws.Range("D54").Value = 'https://shopping-phinf.pstatic.net/main_6512970/6512970130.20220111134842.jpg'
$ws.Range("E54").Value = '''60360'
$ws.Range("I54").Value = '아이클'
$ws.Range("J54").Value = '아이클'
# Row 55
$ws.Range("B55").Value = '벤큐 WiT 아이케어 LED 스탠드'
$ws.Range("C55").Value = 'https://search.shopping.naver.com/gate.nhn?id=9405434840'
$ws.Range("D55").Value = 'https://shopping-phinf.pstatic.net/main_9405434/9405434840.20201215164048.jpg'
$ws.Range("E55").Value = '''198990'
$ws.Range("G55").Value = '네이버'
$ws.Range("H55").Value = '일반 - 가격비교 상품'
$ws.Range("I55").Value = '벤큐'
$ws.Range("J55").Value = '벤큐'
# Row 56
$ws.Range("B56").Value = '파나소닉 LED스탠드 5W 접이식 무선스탠드 휴대용스탠드 USB충전방식 침대독서등'
$ws.Range("C56").Value = 'https://search.shopping.naver.com/gate.nhn?id=82510260293'
$ws.Range("D56").Value = 'https://shopping-phinf.pstatic.net/main_8251026/82510260293.3.jpg'
$ws.Range("E56").Value = '''23500'
$ws.Range("G56").Value = '엔셀라이트'
$ws.Range("H56").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I56").Value = '파나소닉'
$ws.Range("J56").Value = '파나소닉'
# Row 58
$ws.Range("E58").Value = '''68830'
# Row 62
$ws.Range("B62").Value = '아모전자 아모램프 에스2'
$ws.Range("C62").Value = 'https://search.shopping.naver.com/gate.nhn?id=21444465655'
$ws.Range("D62").Value = 'https://shopping-phinf.pstatic.net/main_2144446/21444465655.20191204140327.jpg'
$ws.Range("E62").Value = '''49000'
$ws.Range("G62").Value = '네이버'
$ws.Range("H62").Value = '일반 - 가격비교 상품'
$ws.Range("I62").Value = '아모램프'
$ws.Range("J62").Value = '아모전자'
# Row 63
$ws.Range("B63").Value = '책상스탠드 LED스탠드 공부스탠드 시력보호 학생 공부방 탁상'
$ws.Range("C63").Value = 'https://search.shopping.naver.com/gate.nhn?id=82294931364'
$ws.Range("D63").Value = 'https://shopping-phinf.pstatic.net/main_8229493/82294931364.6.jpg'
$ws.Range("E63").Value = '''78900'
$ws.Range("G63").Value = '최고의선택.'
$ws.Range("H63").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I63").Value = 'TENEY'
$ws.Range("J63").Value = '태인일렉콤'
# Row 65
$ws.Range("B65").Value = 'e스마트 터치 LED스탠드'
$ws.Range("C65").Value = 'https://search.shopping.naver.com/gate.nhn?id=13922662497'
$ws.Range("D65").Value = 'https://shopping-phinf.pstatic.net/main_1392266/13922662497.20181210164958.jpg'
$ws.Range("E65").Value = '''64420'
$ws.Range("G65").Value = '네이버'
$ws.Range("H65").Value = '일반 - 가격비교 상품'
$ws.Range("I65").Value = 'e스마트'
# Row 66
$ws.Range("B66").Value = '클래시 LED 폴딩 블루라이트차단 스탠드 공부 공부용 학생 학습용 무선스탠드 책상스탠드'
$ws.Range("C66").Value = 'https://search.shopping.naver.com/gate.nhn?id=83068920686'
$ws.Range("D66").Value = 'https://shopping-phinf.pstatic.net/main_8306892/83068920686.9.jpg'
$ws.Range("E66").Value = '''61900'
$ws.Range("G66").Value = '클래시스토어'
$ws.Range("H66").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I66").Value = ""
# Row 67
$ws.Range("E67").Value = '''32710'
# Row 75
$ws.Range("E75").Value = '''23310'
# Row 90
$ws.Range("B90").Value = '피티코퍼레이션 단순생활 LED 책상 스탠드'
$ws.Range("C90").Value = 'https://search.shopping.naver.com/gate.nhn?id=27665714522'
$ws.Range("D90").Value = 'https://shopping-phinf.pstatic.net/main_2766571/27665714522.20210713150202.jpg'
$ws.Range("E90").Value = '''35690'
$ws.Range("G90").Value = '네이버'
$ws.Range("H90").Value = '일반 - 가격비교 상품'
$ws.Range("I90").Value = '단순생활'
$ws.Range("J90").Value = '피티코퍼레이션'
# Row 91
$ws.Range("B91").Value = '파나소닉 접이식 LED스탠드 5W 무선스탠드 휴대용스탠드 USB충전방식 침대독서등 신학기'
$ws.Range("C91").Value = 'https://search.shopping.naver.com/gate.nhn?id=82512827486'
$ws.Range("D91").Value = 'https://shopping-phinf.pstatic.net/main_8251282/82512827486.1.jpg'
$ws.Range("E91").Value = '''23500'
$ws.Range("G91").Value = '빛과 바람'
$ws.Range("H91").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I91").Value = '파나소닉'
$ws.Range("J91").Value = '파나소닉'
# Row 92
$ws.Range("B92").Value = '듀플렉스 DP-910LS'
$ws.Range("C92").Value = 'https://search.shopping.naver.com/gate.nhn?id=10173799151'
$ws.Range("D92").Value = 'https://shopping-phinf.pstatic.net/main_1017379/10173799151.20201013143216.jpg'
$ws.Range("E92").Value = '''29450'
$ws.Range("G92").Value = '네이버'
$ws.Range("H92").Value = '일반 - 가격비교 상품'
$ws.Range("I92").Value = '듀플렉스'
$ws.Range("J92").Value = '듀플렉스'
# Row 93
$ws.Range("B93").Value = '3M LED스탠드 Air X+/10 시력보호 독서실 학생용 공부 책상 스탠드 조명'
$ws.Range("C93").Value = 'https://search.shopping.naver.com/gate.nhn?id=82411579763'
$ws.Range("D93").Value = 'https://shopping-phinf.pstatic.net/main_8241157/82411579763.5.jpg'
$ws.Range("E93").Value = '''58990'
$ws.Range("G93").Value = '스페이스작'
$ws.Range("I93").Value = '3M'
$ws.Range("J93").Value = '3M'
# Row 94
$ws.Range("B94").Value = '파파 LED와이드스탠드 800B 500 책상 책상용 학습용 사무용 공부용 조명 독서등 스탠드 스텐드'
$ws.Range("C94").Value = 'https://search.shopping.naver.com/gate.nhn?id=82962775566'
$ws.Range("D94").Value = 'https://shopping-phinf.pstatic.net/main_8296277/82962775566.3.jpg'
$ws.Range("E94").Value = '''40500'
$ws.Range("G94").Value = '주랩'
$ws.Range("H94").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I94").Value = '파파'
$ws.Range("J94").Value = '베스트조명'
# Row 95
$ws.Range("B95").Value = '필립스 데스크 라이트 갈릴레오 66102'
$ws.Range("C95").Value = 'https://search.shopping.naver.com/gate.nhn?id=25748863523'
$ws.Range("D95").Value = 'https://shopping-phinf.pstatic.net/main_2574886/25748863523.20210126175348.jpg'
$ws.Range("E95").Value = '''56000'
$ws.Range("I95").Value = '필립스'
$ws.Range("J95").Value = ""
# Row 96
$ws.Range("B96").Value = '대원씨엘 디트렌드 루미나 앱솔루트'
$ws.Range("C96").Value = 'https://search.shopping.naver.com/gate.nhn?id=25680612522'
$ws.Range("D96").Value = 'https://shopping-phinf.pstatic.net/main_2568061/25680612522.20210217105610.jpg'
$ws.Range("E96").Value = '''64350'
$ws.Range("G96").Value = '네이버'
$ws.Range("H96").Value = '일반 - 가격비교 상품'
$ws.Range("I96").Value = '디트렌드'
$ws.Range("J96").Value = '대원씨엘'
# Row 97
$ws.Range("B97").Value = '플랜룩스 무선스탠드 LSV-01 타이머 책상 LED스탠드 침대 독서등 충전식'
$ws.Range("C97").Value = 'https://search.shopping.naver.com/gate.nhn?id=81371921035'
$ws.Range("D97").Value = 'https://shopping-phinf.pstatic.net/main_8137192/81371921035.4.jpg'
$ws.Range("E97").Value = '''24900'
$ws.Range("G97").Value = '플랜룩스 스토어'
$ws.Range("H97").Value = '일반 - 가격비교 비매칭 일반상품'
$ws.Range("I97").Value = '플랜룩스'
$ws.Range("J97").Value = '플랜룩스'
# Row 98
$ws.Range("B98").Value = '필립스 할리 66127'
$ws.Range("C98").Value = 'https://search.shopping.naver.com/gate.nhn?id=21327507911'
$ws.Range("D98").Value = 'https://shopping-phinf.pstatic.net/main_2132750/21327507911.20191118105932.jpg'
$ws.Range("E98").Value = '''32940'
$ws.Range("J98").Value = '필립스'
# Row 99
$ws.Range("B99").Value = '넥소버 NXL-5000'
$ws.Range("C99").Value = 'https://search.shopping.naver.com/gate.nhn?id=27805503522'
$ws.Range("D99").Value = 'https://shopping-phinf.pstatic.net/main_2780550/27805503522.20211110161237.jpg'
$ws.Range("E99").Value = '''23160'
$ws.Range("I99").Value = '넥소버'
$ws.Range("J99").Value = '넥소버'
